$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "IMDB" column header in F1
$ws.Range("F1").Value = "IMDB"

# Add IMDB rating values for each film row (F2:F9)
$ws.Range("F2").Value = 4.5
$ws.Range("F3").Value = 4.9
$ws.Range("F4").Value = 6.2
$ws.Range("F5").Value = 8.2
$ws.Range("F6").Value = 4.1
$ws.Range("F7").Value = 9.4
$ws.Range("F8").Value = 3.5
$ws.Range("F9").Value = 6.6

# Update the selection to match the post-edit state (F10)
$ws.Range("F10").Select()
